# Auto-generated script to update market price columns (H-N) across all sheets
# per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 200
$ws.Cells.Item(4, 9).Value = 200
$ws.Cells.Item(4, 11).Value = 200
$ws.Cells.Item(4, 13).Value = -86
$ws.Cells.Item(43, 8).Value = 2090.8
$ws.Cells.Item(43, 9).Value = 1775
$ws.Cells.Item(43, 10).Value = 2301.3333
$ws.Cells.Item(43, 11).Value = 1775
$ws.Cells.Item(43, 12).Value = 2301.3333
$ws.Cells.Item(43, 13).Value = -1706
$ws.Cells.Item(43, 14).Value = -2439.3333
$ws.Cells.Item(58, 8).Value = 1049.6
$ws.Cells.Item(58, 9).Value = 86.40000000000001
$ws.Cells.Item(58, 10).Value = 2012.8
$ws.Cells.Item(58, 11).Value = 259.2
$ws.Cells.Item(58, 12).Value = 6038.4
$ws.Cells.Item(58, 13).Value = -109.2
$ws.Cells.Item(58, 14).Value = -6338.4
$ws.Cells.Item(92, 8).Value = 1167.375
$ws.Cells.Item(92, 9).Value = 1167.375
$ws.Cells.Item(92, 11).Value = 1167.375
$ws.Cells.Item(92, 13).Value = 80.625
$ws.Cells.Item(127, 8).Value = 1065.6666
$ws.Cells.Item(127, 9).Value = 1065.6666
$ws.Cells.Item(127, 11).Value = 3196.9998
$ws.Cells.Item(127, 13).Value = 1763.0002
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 13).ClearContents()
$ws.Cells.Item(131, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 8763.75
$ws.Cells.Item(132, 9).Value = 8763.75
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 26291.25
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -23761.25
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 2392.3333
$ws.Cells.Item(137, 9).Value = 2392.3333
$ws.Cells.Item(137, 11).Value = 7176.999899999999
$ws.Cells.Item(137, 13).Value = -4626.999899999999
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 13).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 591
$ws.Cells.Item(26, 9).Value = 591
$ws.Cells.Item(26, 11).Value = 591
$ws.Cells.Item(26, 13).Value = -261
$ws.Cells.Item(32, 8).Value = 9936.888999999999
$ws.Cells.Item(32, 9).Value = 8168.5293
$ws.Cells.Item(32, 11).Value = 8168.5293
$ws.Cells.Item(32, 13).Value = -7881.5293
$ws.Cells.Item(74, 8).Value = 1361
$ws.Cells.Item(74, 9).Value = 1361
$ws.Cells.Item(74, 11).Value = 1361
$ws.Cells.Item(74, 13).Value = -487
$ws.Cells.Item(77, 8).Value = 1361
$ws.Cells.Item(77, 9).Value = 1361
$ws.Cells.Item(77, 11).Value = 6805
$ws.Cells.Item(77, 13).Value = -2437

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 45074
$ws.Cells.Item(35, 10).Value = 45074
$ws.Cells.Item(35, 12).Value = 45074
$ws.Cells.Item(35, 14).Value = -45694
$ws.Cells.Item(99, 8).Value = 1041.4286
$ws.Cells.Item(99, 10).Value = 1099
$ws.Cells.Item(99, 12).Value = 1099
$ws.Cells.Item(99, 14).Value = -4095

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 31.4
$ws.Cells.Item(7, 9).Value = 10.5
$ws.Cells.Item(7, 10).Value = 62.75
$ws.Cells.Item(7, 11).Value = 10.5
$ws.Cells.Item(7, 12).Value = 62.75
$ws.Cells.Item(7, 13).Value = 102.5
$ws.Cells.Item(7, 14).Value = -288.75
$ws.Cells.Item(22, 8).Value = 549.5
$ws.Cells.Item(22, 9).Value = 399
$ws.Cells.Item(22, 11).Value = 399
$ws.Cells.Item(22, 13).Value = -49
$ws.Cells.Item(39, 8).Value = 2666.3333
$ws.Cells.Item(39, 9).Value = 2666.3333
$ws.Cells.Item(39, 11).Value = 2666.3333
$ws.Cells.Item(39, 13).Value = -2275.3333
$ws.Cells.Item(49, 8).Value = 2666.3333
$ws.Cells.Item(49, 9).Value = 2666.3333
$ws.Cells.Item(49, 11).Value = 2666.3333
$ws.Cells.Item(49, 13).Value = -2484.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 30
$ws.Cells.Item(26, 9).Value = 30
$ws.Cells.Item(26, 11).Value = 90
$ws.Cells.Item(26, 13).Value = 198
$ws.Cells.Item(34, 8).Value = 2166
$ws.Cells.Item(34, 10).Value = 2999
$ws.Cells.Item(34, 12).Value = 8997
$ws.Cells.Item(34, 14).Value = -9165
$ws.Cells.Item(39, 8).Value = 1514.1428
$ws.Cells.Item(39, 10).Value = 2999.6667
$ws.Cells.Item(39, 12).Value = 8999.000100000001
$ws.Cells.Item(39, 14).Value = -9587.000100000001
$ws.Cells.Item(55, 8).Value = 2170.6
$ws.Cells.Item(55, 10).Value = 2488.25
$ws.Cells.Item(55, 12).Value = 7464.75
$ws.Cells.Item(55, 14).Value = -7818.75
$ws.Cells.Item(117, 8).Value = 5279.8
$ws.Cells.Item(117, 9).Value = 705.6667
$ws.Cells.Item(117, 11).Value = 2117.0001
$ws.Cells.Item(117, 13).Value = 1324.9999
$ws.Cells.Item(128, 8).Value = 149999
$ws.Cells.Item(128, 9).Value = 149999
$ws.Cells.Item(128, 11).Value = 449997
$ws.Cells.Item(128, 13).Value = -445017

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 2042884
$ws.Cells.Item(3, 10).Value = 2525050
$ws.Cells.Item(3, 12).Value = 2525050
$ws.Cells.Item(3, 14).Value = -2525282
$ws.Cells.Item(10, 8).Value = 1000
$ws.Cells.Item(10, 9).Value = 1000
$ws.Cells.Item(10, 11).Value = 1000
$ws.Cells.Item(10, 13).Value = -831
$ws.Cells.Item(11, 8).Value = 4722655.5
$ws.Cells.Item(11, 9).Value = 4625475
$ws.Cells.Item(11, 11).Value = 4625475
$ws.Cells.Item(11, 13).Value = -4625336
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).ClearContents()
$ws.Cells.Item(14, 14).ClearContents()
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 14).ClearContents()
$ws.Cells.Item(43, 8).Value = 3758.5
$ws.Cells.Item(43, 9).Value = 3758.5
$ws.Cells.Item(43, 11).Value = 3758.5
$ws.Cells.Item(43, 13).Value = -3607.5
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 14).ClearContents()
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 14569
$ws.Cells.Item(80, 9).Value = 2972.5
$ws.Cells.Item(80, 11).Value = 2972.5
$ws.Cells.Item(80, 13).Value = -1974.5
$ws.Cells.Item(83, 8).Value = 14569
$ws.Cells.Item(83, 9).Value = 2972.5
$ws.Cells.Item(83, 11).Value = 14862.5
$ws.Cells.Item(83, 13).Value = -9870.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(22, 8).Value = 197.66667
$ws.Cells.Item(22, 9).Value = 197.66667
$ws.Cells.Item(22, 11).Value = 197.66667
$ws.Cells.Item(22, 13).Value = 97.33332999999999
$ws.Cells.Item(27, 8).Value = 197.66667
$ws.Cells.Item(27, 9).Value = 197.66667
$ws.Cells.Item(27, 11).Value = 197.66667
$ws.Cells.Item(27, 13).Value = -90.66667000000001
$ws.Cells.Item(58, 8).Value = 26625
$ws.Cells.Item(58, 10).Value = 50000
$ws.Cells.Item(58, 12).Value = 50000
$ws.Cells.Item(58, 14).Value = -50520
$ws.Cells.Item(93, 8).Value = 724.5
$ws.Cells.Item(93, 9).Value = 724.5
$ws.Cells.Item(93, 11).Value = 724.5
$ws.Cells.Item(93, 13).Value = 523.5
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 2357
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 13).ClearContents()
$ws.Cells.Item(14, 8).Value = 6002.5
$ws.Cells.Item(14, 9).Value = 6002.5
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 6002.5
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -5834.5
$ws.Cells.Item(14, 14).ClearContents()
$ws.Cells.Item(23, 8).Value = 99.666664
$ws.Cells.Item(23, 9).Value = 99.666664
$ws.Cells.Item(23, 11).Value = 99.666664
$ws.Cells.Item(23, 13).Value = 129.333336
$ws.Cells.Item(96, 8).Value = 3330
$ws.Cells.Item(96, 9).Value = 3250
$ws.Cells.Item(96, 11).Value = 3250
$ws.Cells.Item(96, 13).Value = -1877
$ws.Cells.Item(100, 8).Value = 430
$ws.Cells.Item(100, 9).Value = 312.5
$ws.Cells.Item(100, 11).Value = 625
$ws.Cells.Item(100, 13).Value = -84
$ws.Cells.Item(107, 8).Value = 1411.1
$ws.Cells.Item(107, 9).Value = 850.8
$ws.Cells.Item(107, 11).Value = 2552.4
$ws.Cells.Item(107, 13).Value = -632.3999999999996
